$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 401, shifting existing rows 401-424 down to 402-425.
$ws.Rows.Item(401).Insert()

# Populate the new row 401 with values (copy the static category fields, and set the new data).
$ws.Cells.Item(401, 1).Value = 3
$ws.Cells.Item(401, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(401, 3).Value = "Coquimbo"
$ws.Cells.Item(401, 4).Value = 44746
$ws.Cells.Item(401, 5).Value = 5
$ws.Cells.Item(401, 6).Value = "Fruta"
$ws.Cells.Item(401, 7).Value = 100108
$ws.Cells.Item(401, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(401, 9).Value = 100108002
$ws.Cells.Item(401, 10).Value = "Mango"
$ws.Cells.Item(401, 11).Value = "Sin especificar"
$ws.Cells.Item(401, 12).Value = "Primera"
$ws.Cells.Item(401, 13).Value = 456
$ws.Cells.Item(401, 14).Value = 8500
$ws.Cells.Item(401, 15).Value = 9000
$ws.Cells.Item(401, 16).Value = 8750
$ws.Cells.Item(401, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(401, 18).Value = "Brasil"
$ws.Cells.Item(401, 19).Value = 2188
$ws.Cells.Item(401, 20).Value = 4
